# Add new enum values (MRM, PRM) for reaction monitoring to the
# data_collection_mode list, update the dependent data validation on the
# main sheet, and update the explanatory comment on the "data_collection_mode"
# column header (X1) to mention the new values.

$wb = $excel.ActiveWorkbook

# 1. Add MRM / PRM to the "data_collection_mode list" helper sheet.
$modeList = $wb.Worksheets.Item("data_collection_mode list")
$modeList.Range("A3").Value = "MRM"
$modeList.Range("A4").Value = "PRM"

# 2. Update the data validation on column X (data_collection_mode) of the
#    main "Export as TSV" sheet so it references the expanded list and the
#    error message reflects the new allowed values.
$main = $wb.Worksheets.Item("Export as TSV")
$validation = $main.Range("X2:X1048576").Validation
$validation.Formula1 = "'data_collection_mode list'!`$A`$1:`$A`$4"
$validation.ErrorMessage = "Value must be one of: DDA / DIA / MRM / PRM."

# 3. Update the cell comment on X1 describing the allowed values.
$comment = $main.Range("X1").Comment
$comment.Text("Mode of data collection in tandem MS assays. Either DDA (Data-dependent acquisition), DIA (Data-independent acquisition), MRM (multiple reaction monitoring), or PRM (parallel reaction monitoring).")
